$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.992.24"
$ws.Range("E2").Value = "  -0.40%  "
$ws.Range("D3").Value = "1.859.91"
$ws.Range("E3").Value = "  -0.88%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").Value = "311.95"
$ws.Range("E5").Value = "  -0.50%  "
$ws.Range("D7").Value = "0.5147"
$ws.Range("E7").Value = "  +0.97%  "
$ws.Range("D8").Value = "0.3835"
$ws.Range("E8").Value = "  -0.38%  "
$ws.Range("D9").Value = "0.08215"
$ws.Range("E9").Value = "  -9.84%  "
$ws.Range("D10").Value = "1.110"
$ws.Range("E10").Value = "  -1.18%  "
$ws.Range("D11").Value = "41.51"
$ws.Range("E11").Value = "  -0.16%  "
$ws.Range("D12").Value = "6.191"
$ws.Range("E12").Value = "  -2.46%  "
$ws.Range("D13").Value = "20.56"
$ws.Range("E13").Value = "  -1.02%  "
$ws.Range("D14").Value = "1.861.90"
$ws.Range("E14").Value = "  -0.71%  "
$ws.Range("D15").Value = "7.251"
$ws.Range("E15").Value = "  +0.55%  "
$ws.Range("E16").Value = "  +0.11%  "
$ws.Range("E17").Value = "  -1.63%  "
$ws.Range("D18").Value = "90.59"
$ws.Range("E18").Value = "  -0.59%  "
$ws.Range("D19").Value = "0.06644"
$ws.Range("E19").Value = "  +0.71%  "
$ws.Range("E20").Value = "  -2.63%  "
$ws.Range("D21").Value = "1.002"
$ws.Range("E21").Value = "  +0.07%  "
$ws.Range("D22").Value = "6.007"
$ws.Range("E22").Value = "  -1.71%  "
$ws.Range("D23").Value = "28.015.49"
$ws.Range("E23").Value = "  -0.38%  "
$ws.Range("D24").Value = "11.07"
$ws.Range("E24").Value = "  -3.11%  "
$ws.Range("E25").Value = "  -1.59%  "
$ws.Range("D26").Value = "2.072.63"
$ws.Range("E26").Value = "  -0.89%  "
$ws.Range("D27").Value = "2.510"
$ws.Range("E27").Value = "  -1.24%  "
$ws.Range("D28").Value = "157.92"
$ws.Range("E28").Value = "  +0.16%  "
$ws.Range("D29").Value = "20.45"
$ws.Range("E29").Value = "  -1.75%  "
$ws.Range("D30").Value = "124.71"
$ws.Range("E30").Value = "  -1.61%  "
$ws.Range("E31").Value = "  +0.94%  "
$ws.Range("E32").Value = "  -3.19%  "
$ws.Range("D33").Value = "5.945"
$ws.Range("E33").Value = "  +5.89%  "
$ws.Range("D34").Value = "3.592"
$ws.Range("E34").Value = "  -0.16%  "
$ws.Range("D35").Value = "9.354"
$ws.Range("E35").Value = "  -3.38%  "
$ws.Range("D36").Value = "0.02416"
$ws.Range("E36").Value = "  -0.76%  "
$ws.Range("D37").Value = "0.06502"
$ws.Range("E37").Value = "  -1.08%  "
$ws.Range("D38").Value = "0.2174"
$ws.Range("E38").Value = "  -0.14%  "
$ws.Range("D39").Value = "0.6540"
$ws.Range("E40").Value = "  -1.31%  "
$ws.Range("D41").Value = "4.992"
$ws.Range("E41").Value = "  +1.58%  "
$ws.Range("D42").Value = "1.215"
$ws.Range("E42").Value = "  -3.83%  "
$ws.Range("D43").Value = "11.16"
$ws.Range("E43").Value = "  -3.62%  "
$ws.Range("D44").Value = "0.6154"
$ws.Range("E44").Value = "  +2.27%  "
$ws.Range("D45").Value = "13.05"
$ws.Range("E45").Value = "  -1.47%  "
$ws.Range("D46").Value = "1.283"
$ws.Range("E46").Value = "  +0.62%  "
$ws.Range("D47").Value = "3.668"
$ws.Range("E47").Value = "  -0.08%  "
$ws.Range("D48").Value = "2.003"
$ws.Range("E48").Value = "  +0.21%  "
$ws.Range("D49").Value = "1.217"
$ws.Range("E49").Value = "  -1.67%  "
$ws.Range("D50").Value = "120.59"
$ws.Range("E50").Value = "  -0.56%  "
$ws.Range("D51").Value = "78.13"
$ws.Range("E51").Value = "  -2.01%  "
